# PurchaseList.xlsx update:
#   - revert to old smaller size: drop the extra "basket" line item (Designator
#     "B1", Comment/Footprint "N-5") that lived in row 23 of the BOM sheet.
#   - a few supplier stock numbers were refreshed for the remaining parts.
#   - the "printed:" timestamp moved from 16:31 to 17:43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Supplier Stock 1" numbers (column J) for a handful of rows.
$ws.Range("J2").Value = 1903578
$ws.Range("J7").Value = 43490
$ws.Range("J17").Value = 4399293

# Remove the obsolete "B1" / "N-5" basket row entirely; everything below
# shifts up by one row (row 27 -> row 26, dimension A1:R27 -> A1:R26).
$ws.Rows("23").Delete()

# Bump the "printed:" time label next to the totals row (now row 24 after
# the deletion above). Keep it a literal text value (quote-prefixed) so it
# is stored the same way as the original "16:31" label.
$ws.Range("F24").Value = "'17:43"
